# Rename the worksheet to reflect the unified "DataNode" concept
# (Property1 -> DataNode), per the commit:
#   "unify the conception of DataNode, DataTable, Entity."
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"

# Preserve the author's last on-screen selection in the frozen pane
# (moved from A9 to D40) when the workbook was saved.
$ws.Range("D40").Select()
